$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "Poor smartphone photo scans are really annoying and these researchers finally fi"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2714"

$ws.Range("D32").Value = "한개의 모델로 성격이 비슷한 여러개의 모델을 대체해보자"
$ws.Range("E32").Value = "https://dodonam.tistory.com/299"

$ws.Range("D39").Value = "Facial Expression Recognition with Keras"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Facial-Expression-Recognition-with-Keras-1"

$ws.Range("D46").Value = "[한국수력원자력] 2021년 02월, 생물정보학(Bioinformatics 채용), 보건의료 빅데이터 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/378"
